# "[Refactor] add final edits to manuscript"
#
# 1. Title (Heading1) paragraph: line spacing changes from double (w:line="480")
#    to single (w:line="240"), both w:lineRule="auto".
# 2. Abstract (Heading2) paragraph: gains an explicit single line-spacing
#    (w:line="240" w:lineRule="auto") it did not have before.
# 3. The paragraph that previously held only the page-break run (and an
#    otherwise-empty w:spacing w:line="480" pPr) is split in two:
#      - a new "Main text" paragraph styled BodyText with
#        w:lang w:eastAsia="de-DE" on both the paragraph mark and the run
#      - a following paragraph that keeps just the page-break run, now with
#        no paragraph properties at all.

$d = $word.ActiveDocument

# --- 1) Title: double -> single line spacing -------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.LineSpacingRule = 0   # wdLineSpaceSingle -> <w:spacing w:line="240" w:lineRule="auto"/>

# --- 2) Abstract heading: add single line spacing --------------------------
$abstractPara = $d.Paragraphs.Item(2)
$abstractPara.LineSpacingRule = 0   # adds <w:spacing w:line="240" w:lineRule="auto"/>

# --- 3) Split the page-break paragraph into "Main text" + page break -------
$breakPara = $d.Paragraphs.Item(3)
$breakRange = $d.Range($breakPara.Range.Start, $breakPara.Range.End)

$packageXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="BodyText"/>
              <w:rPr>
                <w:lang w:eastAsia="de-DE"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:lang w:eastAsia="de-DE"/>
              </w:rPr>
              <w:t>Main text</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:br w:type="page"/>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$breakRange.InsertXML($packageXml)

Write-Host "Applied: line-spacing tweaks + inserted 'Main text' paragraph before page break."
